# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on the active worksheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.558.63'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.920.73'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.36'
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4817'
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4059'
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.011'
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.42'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.927.48'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.065'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.246'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.69'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06867'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.60'
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.560.30'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.684'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.90'
$ws.Range("E23").Value = '  +0.86%  '
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.139.52'
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.536'
$ws.Range("E26").Value = '  +3.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.03'
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.00'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.098'
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.72'
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.020'
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09634'
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.631'
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.558'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.378'
$ws.Range("E35").Value = '  -1.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06351'
$ws.Range("E36").Value = '  +3.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02288'
$ws.Range("E37").Value = '  +0.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.200'
$ws.Range("E38").Value = '  +2.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5950'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.74'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.922'
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1851'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.496'
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("E44").Value = '  +2.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.43'
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07495'
$ws.Range("E46").Value = '  -2.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5572'
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.945'
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.63'
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.433'
$ws.Range("E50").Value = '  +3.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.22'
$ws.Range("E51").Value = '  -0.89%  '
